# Apply updated cryptos list values (generated from commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '54.541.52'
$ws.Range("E2").Value = '  +0.90%  '

# Row 3
$ws.Range("D3").Value = '2.288.55'
$ws.Range("E3").Value = '  +0.11%  '

# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.39%  '

# Row 5
$ws.Range("D5").Value = "'503.37"
$ws.Range("E5").Value = '  +1.91%  '

# Row 6
$ws.Range("D6").Value = "'130.48"
$ws.Range("E6").Value = '  +2.76%  '

# Row 7
$ws.Range("E7").Value = '  -0.23%  '

# Row 8
$ws.Range("E8").Value = '  +0.57%  '

# Row 9
$ws.Range("E9").Value = '  +1.73%  '

# Row 10
$ws.Range("E10").Value = '  +0.68%  '

# Row 11
$ws.Range("D11").Value = "'0.338"
$ws.Range("E11").Value = '  +4.97%  '

# Row 12
$ws.Range("E12").Value = '  +2.76%  '

# Row 13
$ws.Range("D13").Value = '2.699.90'
$ws.Range("E13").Value = '  +0.25%  '

# Row 14
$ws.Range("D14").Value = "'22.90"
$ws.Range("E14").Value = '  +6.43%  '

# Row 15
$ws.Range("D15").Value = '54.438.69'
$ws.Range("E15").Value = '  +0.38%  '

# Row 16
$ws.Range("E16").Value = '  +0.70%  '

# Row 17
$ws.Range("D17").Value = '2.309.40'
$ws.Range("E17").Value = '  +0.02%  '

# Row 18
$ws.Range("D18").Value = "'10.29"
$ws.Range("E18").Value = '  +3.30%  '

# Row 19
$ws.Range("E19").Value = '  +3.02%  '

# Row 20
$ws.Range("D20").Value = "'304.59"
$ws.Range("E20").Value = '  +0.60%  '

# Row 21
$ws.Range("D21").Value = "'6.37"
$ws.Range("E21").Value = '  -0.69%  '

# Row 22
$ws.Range("E22").Value = '  -0.12%  '

# Row 23
$ws.Range("D23").Value = "'61.90"
$ws.Range("E23").Value = '  -2.72%  '

# Row 24
$ws.Range("E24").Value = '  -0.48%  '

# Row 25
$ws.Range("E25").Value = '  +1.95%  '

# Row 26
$ws.Range("E26").Value = '  +3.61%  '

# Row 27
$ws.Range("D27").Value = "'171.72"
$ws.Range("E27").Value = '  +2.09%  '

# Row 28
$ws.Range("E28").Value = '  +2.19%  '

# Row 29
$ws.Range("D29").Value = '0.0₃0695'
$ws.Range("E29").Value = '  +1.60%  '

# Row 30
$ws.Range("E30").Value = '  +1.69%  '

# Row 31
$ws.Range("E31").Value = '  +0.84%  '

# Row 32
$ws.Range("E32").Value = '  -0.05%  '

# Row 33
$ws.Range("D33").Value = "'17.85"
$ws.Range("E33").Value = '  +1.51%  '

# Row 34
$ws.Range("D34").Value = "'0.970"
$ws.Range("E34").Value = '  +11.02%  '

# Row 35
$ws.Range("D35").Value = "'0.997"
$ws.Range("E35").Value = '  +0.12%  '

# Row 36
$ws.Range("E36").Value = '  +0.88%  '

# Row 37
$ws.Range("D37").Value = "'3.74"
$ws.Range("E37").Value = '  +3.35%  '

# Row 38
$ws.Range("E38").Value = '  +0.49%  '

# Row 39
$ws.Range("E39").Value = '  +1.42%  '

# Row 40
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").Value = "'5.07"
$ws.Range("E40").Value = '  +5.41%  '

# Row 41
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = "'3.39"
$ws.Range("E41").Value = '  +1.58%  '

# Row 42
$ws.Range("D42").Value = "'126.37"
$ws.Range("E42").Value = '  -0.36%  '

# Row 43
$ws.Range("E43").Value = '  +3.56%  '

# Row 44
$ws.Range("E44").Value = '  +1.03%  '

# Row 46
$ws.Range("D46").Value = "'242.18"
$ws.Range("E46").Value = '  +1.24%  '

# Row 47
$ws.Range("E47").Value = '  +0.51%  '

# Row 48
$ws.Range("E48").Value = '  +1.80%  '

# Row 49
$ws.Range("E49").Value = '  +0.78%  '

# Row 50
$ws.Range("D50").Value = "'16.45"
$ws.Range("E50").Value = '  +1.11%  '

# Row 51
$ws.Range("E51").Value = '  -0.57%  '
